$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly data points (Aug-22 through Dec-22), continuing the feed.
# Column A: date (serial), formatted like the existing rows (mmm-yy),
# alternating with the existing "highlight" fill used on the previous
# last row (176) every other row.
# Column B: netback value.

$rows = @(
    @{ Row = 177; Date = 44774; Value = 748.22138584332561; Highlight = $false },
    @{ Row = 178; Date = 44805; Value = 752.12325771793871; Highlight = $true  },
    @{ Row = 179; Date = 44835; Value = 753.07239912746263; Highlight = $false },
    @{ Row = 180; Date = 44866; Value = 752.53385046079586; Highlight = $true  },
    @{ Row = 181; Date = 44896; Value = 751.62560635603404; Highlight = $false }
)

foreach ($r in $rows) {
    $aCell = $ws.Cells.Item($r.Row, 1)
    $bCell = $ws.Cells.Item($r.Row, 2)

    $aCell.Value = $r.Date
    $aCell.NumberFormat = "mmm-yy"
    if ($r.Highlight) {
        $aCell.Interior.Color = 65535
    }

    $bCell.Value = $r.Value
}

# Match the saved selection/view state from the commit.
$null = $ws.Range("A176").Select()
